$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("openbis-metadata")
$ws2 = $wb.Worksheets.Item("openbis-data")

# Remove the "Strain" row (row 3) from the metadata sheet; this shifts all
# subsequent rows up by one and Excel will drop now-unused shared strings.
$ws1.Rows.Item(3).Delete()

# The "openbis-data" sheet's example header value changes from "Abs" to "Strain".
$ws2.Range("A1").Value = "Strain"

# Restore the selections/active cells to match the authored state.
$ws2.Activate()
$ws2.Range("A2").Select()

$ws1.Activate()
$ws1.Range("A11").Select()
